$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 118.85714
$ws.Range("I5").Value = 118.85714
$ws.Range("K5").Value = 118.85714
$ws.Range("M5").Value = -3.857140000000001
$ws.Range("H28").Value = 268.25
$ws.Range("I28").Value = 155.58333
$ws.Range("J28").Value = 606.25
$ws.Range("K28").Value = 155.58333
$ws.Range("L28").Value = 606.25
$ws.Range("M28").Value = 329.41667
$ws.Range("N28").Value = -1576.25
$ws.Range("H40").Value = 1899.9166
$ws.Range("I40").Value = 1879.9
$ws.Range("K40").Value = 1879.9
$ws.Range("M40").Value = -1704.9
$ws.Range("H43").Value = 5088.8
$ws.Range("I43").Value = 5814.6665
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 5814.6665
$ws.Range("L43").Value = 4000
$ws.Range("M43").Value = -5745.6665
$ws.Range("N43").Value = -4138
$ws.Range("H86").Value = 3720.3
$ws.Range("I86").Value = 2695.8
$ws.Range("K86").Value = 2695.8
$ws.Range("M86").Value = -1572.8
$ws.Range("H89").Value = 3720.3
$ws.Range("I89").Value = 2695.8
$ws.Range("K89").Value = 13479
$ws.Range("M89").Value = -7863
$ws.Range("H99").Value = 204.33333
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7400
$ws.Range("I31").Value = 7400
$ws.Range("K31").Value = 7400
$ws.Range("M31").Value = -7106
$ws.Range("H32").Value = 24735.484
$ws.Range("I32").Value = 20489.055
$ws.Range("K32").Value = 20489.055
$ws.Range("M32").Value = -20202.055
$ws.Range("H61").Value = 1600
$ws.Range("I61").Value = 900
$ws.Range("K61").Value = 900
$ws.Range("M61").Value = -688
$ws.Range("H112").Value = 100000
$ws.Range("J112").Value = 100000
$ws.Range("L112").Value = 100000
$ws.Range("N112").Value = -102954
$ws.Range("H114").Value = 59999
$ws.Range("J114").Value = 59999
$ws.Range("L114").Value = 59999
$ws.Range("N114").Value = -68677
$ws.Range("H132").Value = 2294.077
$ws.Range("I132").Value = 2252.4
$ws.Range("J132").Value = 2433
$ws.Range("K132").Value = 6757.200000000001
$ws.Range("L132").Value = 7299
$ws.Range("M132").Value = -4227.200000000001
$ws.Range("N132").Value = -12359
$ws.Range("H136").Value = 1600
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = -150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8279.333000000001
$ws.Range("J86").Value = 7919.25
$ws.Range("L86").Value = 7919.25
$ws.Range("N86").Value = -10165.25
$ws.Range("H89").Value = 8279.333000000001
$ws.Range("J89").Value = 7919.25
$ws.Range("L89").Value = 39596.25
$ws.Range("N89").Value = -50828.25
$ws.Range("H102").Value = 9000
$ws.Range("I102").Value = 9000
$ws.Range("K102").Value = 9000
$ws.Range("M102").Value = -5755
$ws.Range("H134").Value = 5300
$ws.Range("J134").Value = 5300
$ws.Range("L134").Value = 15900
$ws.Range("N134").Value = -20970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 322.25
$ws.Range("J7").Value = 215.6
$ws.Range("L7").Value = 215.6
$ws.Range("N7").Value = -441.6
$ws.Range("H31").Value = 6816.143
$ws.Range("J31").Value = 6969
$ws.Range("L31").Value = 6969
$ws.Range("N31").Value = -7559
$ws.Range("H34").Value = 6816.143
$ws.Range("J34").Value = 6969
$ws.Range("L34").Value = 6969
$ws.Range("N34").Value = -7373
$ws.Range("H74").Value = 69999
$ws.Range("J74").Value = 69999
$ws.Range("L74").Value = 69999
$ws.Range("N74").Value = -71747
$ws.Range("H77").Value = 69999
$ws.Range("J77").Value = 69999
$ws.Range("L77").Value = 209997
$ws.Range("N77").Value = -218733
$ws.Range("H86").Value = 6176.7417
$ws.Range("I86").Value = 3561.8096
$ws.Range("K86").Value = 3561.8096
$ws.Range("M86").Value = -2438.8096
$ws.Range("H89").Value = 6176.7417
$ws.Range("I89").Value = 3561.8096
$ws.Range("K89").Value = 17809.048
$ws.Range("M89").Value = -12193.048
$ws.Range("H94").Value = 1086.4
$ws.Range("I94").Value = 828
$ws.Range("J94").Value = 1197.1428
$ws.Range("K94").Value = 828
$ws.Range("L94").Value = 1197.1428
$ws.Range("M94").Value = -377
$ws.Range("N94").Value = -2099.1428
$ws.Range("H105").Value = 4831.6665
$ws.Range("J105").Value = 4997.5
$ws.Range("L105").Value = 4997.5
$ws.Range("N105").Value = -8491.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5999.6665
$ws.Range("I113").Value = 4499
$ws.Range("K113").Value = 13497
$ws.Range("M113").Value = -11327
$ws.Range("H114").Value = 746.61536
$ws.Range("I114").Value = 683.3333
$ws.Range("J114").Value = 765.6
$ws.Range("K114").Value = 2049.9999
$ws.Range("L114").Value = 2296.8
$ws.Range("M114").Value = 1204.0001
$ws.Range("N114").Value = -8804.799999999999
$ws.Range("H117").Value = 1880.1111
$ws.Range("J117").Value = 4197.5
$ws.Range("L117").Value = 12592.5
$ws.Range("N117").Value = -19476.5
$ws.Range("H121").Value = 1598.6364
$ws.Range("J121").Value = 2962.5
$ws.Range("L121").Value = 8887.5
$ws.Range("N121").Value = -11507.5
$ws.Range("H129").Value = 1020
$ws.Range("J129").Value = 1133.3334
$ws.Range("L129").Value = 3400.0002
$ws.Range("N129").Value = -13400.0002
$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 10000
$ws.Range("K137").Value = 30000
$ws.Range("M137").Value = -24900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 263.53333
$ws.Range("I2").Value = 41.333332
$ws.Range("K2").Value = 41.333332
$ws.Range("M2").Value = 71.666668
$ws.Range("H18").Value = 40000
$ws.Range("J18").Value = 40000
$ws.Range("L18").Value = 40000
$ws.Range("N18").Value = -40586
$ws.Range("H29").Value = 18999.916
$ws.Range("I29").Value = 19000
$ws.Range("J29").Value = 18999.889
$ws.Range("K29").Value = 19000
$ws.Range("L29").Value = 18999.889
$ws.Range("M29").Value = -18710
$ws.Range("N29").Value = -19579.889
$ws.Range("H43").Value = 5743.25
$ws.Range("J43").Value = 9862
$ws.Range("L43").Value = 9862
$ws.Range("N43").Value = -10164
$ws.Range("H80").Value = 8874.75
$ws.Range("I80").Value = 7750
$ws.Range("K80").Value = 7750
$ws.Range("M80").Value = -6752
$ws.Range("H83").Value = 8874.75
$ws.Range("I83").Value = 7750
$ws.Range("K83").Value = 38750
$ws.Range("M83").Value = -33758
$ws.Range("H107").Value = 223.6
$ws.Range("I107").Value = 204.5
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 204.5
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1715.5
$ws.Range("N107").Value = -4140
$ws.Range("H122").Value = 411656.6
$ws.Range("I122").Value = 75370.92999999999
$ws.Range("J122").Value = 773810.4
$ws.Range("K122").Value = 226112.79
$ws.Range("L122").Value = 2321431.2
$ws.Range("M122").Value = -223662.79
$ws.Range("N122").Value = -2326331.2
$ws.Range("H132").Value = 4947.1763
$ws.Range("I132").Value = 4230.3
$ws.Range("J132").Value = 5971.2856
$ws.Range("K132").Value = 12690.9
$ws.Range("L132").Value = 17913.8568
$ws.Range("M132").Value = -10160.9
$ws.Range("N132").Value = -22973.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4336.727
$ws.Range("I132").Value = 4336.727
$ws.Range("K132").Value = 13010.181
$ws.Range("M132").Value = -10480.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 666.3333
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 500
$ws.Range("M29").Value = -210
$ws.Range("H76").Value = 55950
$ws.Range("I76").Value = 55950
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 55950
$ws.Range("M76").Value = -55635
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 55950
$ws.Range("I79").Value = 55950
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 55950
$ws.Range("M79").Value = -54858
$ws.Range("N79").ClearContents()
$ws.Range("H126").Value = 250977.25
$ws.Range("I126").Value = 333968
$ws.Range("K126").Value = 1001904
$ws.Range("M126").Value = -999434
$ws.Range("H136").Value = 71112.92999999999
$ws.Range("I136").Value = 3188.5557
$ws.Range("J136").Value = 172999.5
$ws.Range("K136").Value = 9565.667099999999
$ws.Range("L136").Value = 518998.5
$ws.Range("M136").Value = -7015.667099999999
$ws.Range("N136").Value = -524098.5
